$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 55.92857
$ws.Range("I11").Value = 55.92857
$ws.Range("K11").Value = 55.92857
$ws.Range("M11").Value = 84.07142999999999

$ws.Range("H18").Value = 7671.2856
$ws.Range("I18").Value = 7671.2856
$ws.Range("K18").Value = 7671.2856
$ws.Range("M18").Value = -7387.2856

$ws.Range("H33").Value = 10883.134
$ws.Range("I33").Value = 11517.643
$ws.Range("K33").Value = 11517.643
$ws.Range("M33").Value = -11288.643

$ws.Range("H40").Value = 5317.4287
$ws.Range("I40").Value = 3502.1
$ws.Range("J40").Value = 6325.9443
$ws.Range("K40").Value = 3502.1
$ws.Range("L40").Value = 6325.9443
$ws.Range("M40").Value = -3327.1
$ws.Range("N40").Value = -6675.9443

$ws.Range("H87").Value = 45000
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 45000
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H131").Value = 42989.348
$ws.Range("I131").Value = 48896.637
$ws.Range("K131").Value = 146689.911
$ws.Range("M131").Value = -141649.911

$ws.Range("H132").Value = 1326.5217
$ws.Range("I132").Value = 1058.6842
$ws.Range("K132").Value = 3176.0526
$ws.Range("M132").Value = -646.0526

$ws.Range("H137").Value = 1772.091
$ws.Range("I137").Value = 1749.3
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 5247.9
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -2697.9
$ws.Range("N137").Value = -11100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3044.5781
$ws.Range("I32").Value = 2535.9465
$ws.Range("K32").Value = 2535.9465
$ws.Range("M32").Value = -2248.9465

$ws.Range("H41").Value = 110683
$ws.Range("I41").Value = 3049
$ws.Range("K41").Value = 3049
$ws.Range("M41").Value = -2635

$ws.Range("H45").Value = 3523.1785
$ws.Range("J45").Value = 6750
$ws.Range("L45").Value = 6750
$ws.Range("N45").Value = -7504

$ws.Range("H97").Value = 604.36
$ws.Range("I97").Value = 565.5294
$ws.Range("J97").Value = 686.875
$ws.Range("K97").Value = 565.5294
$ws.Range("L97").Value = 686.875
$ws.Range("M97").Value = -69.52940000000001
$ws.Range("N97").Value = -1678.875

$ws.Range("H132").Value = 3549.04
$ws.Range("I132").Value = 2761
$ws.Range("K132").Value = 8283
$ws.Range("M132").Value = -5753

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3618
$ws.Range("I22").Value = 4019
$ws.Range("J22").Value = 2615.5
$ws.Range("K22").Value = 4019
$ws.Range("L22").Value = 2615.5
$ws.Range("M22").Value = -3846
$ws.Range("N22").Value = -2961.5

$ws.Range("H86").Value = 3978.5715
$ws.Range("I86").Value = 3978.5715
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3978.5715
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -2855.5715
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 3978.5715
$ws.Range("I89").Value = 3978.5715
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 19892.8575
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -14276.8575
$ws.Range("N89").ClearContents()

$ws.Range("H94").Value = 2648.8096
$ws.Range("I94").Value = 1860.4706
$ws.Range("J94").Value = 5999.25
$ws.Range("K94").Value = 1860.4706
$ws.Range("L94").Value = 5999.25
$ws.Range("M94").Value = -1409.4706
$ws.Range("N94").Value = -6901.25

$ws.Range("H96").Value = 9201
$ws.Range("I96").Value = 9201
$ws.Range("K96").Value = 9201
$ws.Range("M96").Value = -6455

$ws.Range("H107").Value = 3755.913
$ws.Range("I107").Value = 4034.647
$ws.Range("J107").Value = 2966.1667
$ws.Range("K107").Value = 4034.647
$ws.Range("L107").Value = 2966.1667
$ws.Range("M107").Value = -2114.647
$ws.Range("N107").Value = -6806.1667

$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H134").Value = 5401.5415
$ws.Range("I134").Value = 5299.4346
$ws.Range("K134").Value = 15898.3038
$ws.Range("M134").Value = -13363.3038

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4617.3
$ws.Range("J31").Value = 6805.25
$ws.Range("L31").Value = 6805.25
$ws.Range("N31").Value = -7395.25

$ws.Range("H34").Value = 4617.3
$ws.Range("J34").Value = 6805.25
$ws.Range("L34").Value = 6805.25
$ws.Range("N34").Value = -7209.25

$ws.Range("H86").Value = 5984.048
$ws.Range("J86").Value = 7092.7144
$ws.Range("L86").Value = 7092.7144
$ws.Range("N86").Value = -9338.714400000001

$ws.Range("H89").Value = 5984.048
$ws.Range("J89").Value = 7092.7144
$ws.Range("L89").Value = 35463.572
$ws.Range("N89").Value = -46695.572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 350
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 350
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 1050
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -1510

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2284.182
$ws.Range("I132").Value = 2125.111
$ws.Range("K132").Value = 6375.333
$ws.Range("M132").Value = -3845.333

$ws.Range("H141").Value = 68000
$ws.Range("J141").Value = 68000
$ws.Range("L141").Value = 68000
$ws.Range("N141").Value = -78360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 14998.75
$ws.Range("J20").Value = 14998.75
$ws.Range("L20").Value = 14998.75
$ws.Range("N20").Value = -15450.75

$ws.Range("J22").Value = 3000
$ws.Range("L22").Value = 3000
$ws.Range("N22").Value = -3590

$ws.Range("J27").Value = 3000
$ws.Range("L27").Value = 3000
$ws.Range("N27").Value = -3214

$ws.Range("H46").Value = 2717.8
$ws.Range("J46").Value = 3457.6667
$ws.Range("L46").Value = 3457.6667
$ws.Range("N46").Value = -3833.6667

$ws.Range("H93").Value = 11626.238
$ws.Range("I93").Value = 1891.7059
$ws.Range("K93").Value = 1891.7059
$ws.Range("M93").Value = -643.7058999999999

$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").ClearContents()

$ws.Range("H108").Value = 50625
$ws.Range("J108").Value = 50625
$ws.Range("L108").Value = 50625
$ws.Range("N108").Value = -58305

$ws.Range("H136").Value = 2915.7827
$ws.Range("I136").Value = 2717.762
$ws.Range("K136").Value = 8153.286
$ws.Range("M136").Value = -5603.286

$ws.Range("H139").Value = 82907.5
$ws.Range("J139").Value = 82907.5
$ws.Range("L139").Value = 82907.5
$ws.Range("N139").Value = -93187.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 14916.667
$ws.Range("J41").Value = 14916.667
$ws.Range("L41").Value = 14916.667
$ws.Range("N41").Value = -15696.667

$ws.Range("H113").Value = 999.875
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 999.8
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 2999.4
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -7339.4

$ws.Range("H136").Value = 3295.0815
$ws.Range("I136").Value = 2487.3948
$ws.Range("K136").Value = 7462.1844
$ws.Range("M136").Value = -4912.1844

$ws.Range("H140").Value = 97388.5
$ws.Range("J140").Value = 97388.5
$ws.Range("L140").Value = 97388.5
$ws.Range("N140").Value = -107748.5
